$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("E1").Value = "Product"

# Update data row with new values
$ws.Range("A2").Value = "Sanda"
$ws.Range("B2").Value = "Ortiz"
$ws.Range("C2").Value = "micah.littel@gmail.com"
$ws.Range("D2").Value = "fphw2i5ypwj"
$ws.Range("E2").Value = "14.1-inch Laptop"
